# Added range view in Districts
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("District")

# Replace the placeholder "?????" text values in column H (View Range)
# with the actual numeric range values, row by row (row 2 .. row 18).
$values = @(5, 3, 4, 5, 3, 4, 5, 3, 4, 5, 5, 4, 4, 5, 5, 5, 5)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

# Update the active selection to reflect where the user ended up.
$ws.Range("H19").Select()
